$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 1.47
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.8
$ws.Range("AA2").Value = 21
$ws.Range("AC2").Value = 7.5
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 9
$ws.Range("AP2").Value = 26
$ws.Range("AT2").Value = 2.5
$ws.Range("AV2").Value = 67
$ws.Range("BA2").Value = 101

$ws.Range("P5").Value = 3.74

$ws.Range("G6").Value = 2.9
$ws.Range("H6").Value = 2.88
$ws.Range("I6").Value = 2.7
$ws.Range("K6").Value = 2
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.35
$ws.Range("R6").Value = 1.57
$ws.Range("S6").Value = 1.5
$ws.Range("T6").Value = 2.5
$ws.Range("Y6").Value = 11
$ws.Range("AC6").Value = 7
$ws.Range("AG6").Value = 351
$ws.Range("AH6").Value = 7.5
$ws.Range("AT6").Value = 2.5

$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.19

$ws.Range("M8").Value = 1.03
$ws.Range("O8").Value = 1.19
$ws.Range("AX8").Value = 19

$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.22

$ws.Range("M10").Value = 1.02
$ws.Range("O10").Value = 1.15

$ws.Range("G14").Value = 6.5
$ws.Range("L14").Value = 2.25
$ws.Range("M14").Value = 1.1
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 2.63
$ws.Range("Q14").Value = 2.35
$ws.Range("R14").Value = 1.57
$ws.Range("W14").Value = 13
$ws.Range("AA14").Value = 51
$ws.Range("AE14").Value = 23
$ws.Range("AK14").Value = 11
$ws.Range("AL14").Value = 17
$ws.Range("AN14").Value = 7.5
$ws.Range("AX14").Value = 8.5
$ws.Range("AY14").Value = 26
$ws.Range("AZ14").Value = 29

$ws.Range("S15").Value = 1.57

$ws.Range("S16").Value = 1.5

$ws.Range("G17").Value = 1.87
$ws.Range("I17").Value = 3.9
$ws.Range("K17").Value = 2.1
$ws.Range("L17").Value = 4.5
$ws.Range("S17").Value = 1.4
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.91
$ws.Range("X17").Value = 8.5
$ws.Range("AB17").Value = 29
$ws.Range("AC17").Value = 9.5
$ws.Range("AE17").Value = 17
$ws.Range("AG17").Value = 301
$ws.Range("AH17").Value = 10
$ws.Range("AZ17").Value = 81

$ws.Range("G18").Value = 1.92
$ws.Range("S18").Value = 1.44
$ws.Range("T18").Value = 2.63

$ws.Range("I19").Value = 3.4
$ws.Range("M19").Value = 1.08
$ws.Range("N19").Value = 8
